# Generate Report for Handoff
# - Bumps the "Latest Handoff Datetime" for the zh-cn and de-de locales
#   (rows 4-7, files still "Ready for handoff") to reflect a fresh
#   handoff xliff generation pass, and promotes their Priority from
#   "low" to "ht" now that they've been handed off.
# - Propagates the new de-de handoff timestamp up into the Overview
#   sheet's "Latest HO Xliff Generate Date" column for the same rows.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsOverview = $wb.Worksheets.Item("Overview")

# zh-cn: rows 4-7 -> Priority "low" -> "ht"; Latest Handoff Datetime refreshed
$wsZhCn.Range("E4:E7").Value = "ht"
$wsZhCn.Range("H4:H7").Value = "2016-09-01 00:33:45"

# de-de: rows 4-7 -> Priority "low" -> "ht"; Latest Handoff Datetime refreshed
$wsDeDe.Range("E4:E7").Value = "ht"
$wsDeDe.Range("H4:H7").Value = "2016-09-01 00:33:49"

# Overview: rows 4-7 -> Latest HO Xliff Generate Date refreshed to match de-de
$wsOverview.Range("G4:G7").Value = "2016-09-01 00:33:49"
